$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column C ("Förändrad") holds a date serial that was refreshed from
# 45186 (2023-09-17) to 45188 (2023-09-19) for every data row (rows 2-393).
$ws.Range("C2:C393").Value = 45188
